$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "switch" table of days of the week in column F (F3:F10)
$ws.Range("F3").Value = "0 1"
$ws.Range("F4").Value = "2 3"
$ws.Range("F5").Value = "4 5"
$ws.Range("F6").Value = "6 7"
$ws.Range("F7").Value = "8 9"
$ws.Range("F8").Value = 10
$ws.Range("F9").Value = 11
$ws.Range("F10").Value = 12

# Center-align the new values
$ws.Range("F3:F10").HorizontalAlignment = -4108

# Highlight cell C5 in yellow
$ws.Range("C5").Interior.Color = 65535

# Update selection to F4
$ws.Range("F4").Select()

$wb.Save()
